$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells, copying the style used by the existing header row (e.g. AC1)
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record values for each data row (rows 2 through 51)
$lastRow = 51
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 82  # AD
    $ws.Cells.Item($r, 31).Value = 80  # AE
    $ws.Cells.Item($r, 32).Value = 0   # AF
}
